# "solved eqn for weight and true prob from recency and sample size"
#
# Each worksheet holds a probability table keyed by "val" (row A/B column),
# where row 2 is the val=0 edge case. Row 2's over/under split for pts
# (C/D), reb (E/F) and ast (G/H) had been left as a hard-coded 100/0,
# instead of the value implied by solving the weighting equation. The
# corrected figure for val=0 is the mirror (over<->under swapped) of the
# val=1 row (row 3) for that stat - this patch pushes those solved values
# into row 2 for every sheet/stat pair where it had been wrong.
#
# 3pm (I/J) already held the correct value and is left untouched, as are
# the two "8 after" sheets (only 1 data row, nothing to mirror from).

$wb = $excel.ActiveWorkbook

function Set-Row2 {
    param(
        [string]$SheetName,
        [hashtable]$Values
    )
    $ws = $wb.Worksheets.Item($SheetName)
    foreach ($col in $Values.Keys) {
        $ws.Range("${col}2").Value = $Values[$col]
    }
}

Set-Row2 "all 2023 postseason"       @{ C=40; D=60; E=20; F=80; G=80; H=20 }
Set-Row2 "all 2023 full"             @{ C=40; D=60; E=20; F=80; G=80; H=20 }

Set-Row2 "home 2023 postseason"      @{ C=50; D=50; E=25; F=75 }
Set-Row2 "home 2023 full"            @{ C=50; D=50; E=25; F=75 }

Set-Row2 "away 2023 postseason"      @{ C=0;  D=100; E=0;  F=100; G=0;  H=100 }
Set-Row2 "away 2023 full"            @{ C=0;  D=100; E=0;  F=100; G=0;  H=100 }

Set-Row2 "den 2023 postseason"       @{ C=50; D=50; E=50; F=50 }
Set-Row2 "den 2023 full"             @{ C=50; D=50; E=50; F=50 }

Set-Row2 "0 before 2023 postseason"  @{ C=0;  D=100; E=0;  F=100 }
Set-Row2 "0 before 2023 full"        @{ C=0;  D=100; E=0;  F=100 }

Set-Row2 "2 after 2023 postseason"   @{ C=0;  D=100; E=0;  F=100 }
Set-Row2 "2 after 2023 full"         @{ C=0;  D=100; E=0;  F=100 }

Set-Row2 "2 before 2023 postseason"  @{ C=50; D=50; E=50; F=50; G=50; H=50 }
Set-Row2 "2 before 2023 full"        @{ C=50; D=50; E=50; F=50; G=50; H=50 }

Set-Row2 "gsw 2023 postseason"       @{ C=33; D=67; E=0;  F=100; G=67; H=33 }
Set-Row2 "gsw 2023 full"             @{ C=33; D=67; E=0;  F=100; G=67; H=33 }

Set-Row2 "8 before 2023 postseason"  @{ E=0;  F=100 }
Set-Row2 "8 before 2023 full"        @{ E=0;  F=100 }

Set-Row2 "6 after 2023 postseason"   @{ E=0;  F=100 }
Set-Row2 "6 after 2023 full"         @{ E=0;  F=100 }

Set-Row2 "6 before 2023 postseason"  @{ C=0;  D=100; E=0;  F=100 }
Set-Row2 "6 before 2023 full"        @{ C=0;  D=100; E=0;  F=100 }
